$wb = $excel.ActiveWorkbook

# Sheet 1: publico
$ws1 = $wb.Worksheets.Item("publico")
$ws1.Range("D2").Value = 70
$ws1.Range("E2").Value = 50
$ws1.Range("C3").Value = ""
$ws1.Range("E3").Value = 80

# Sheet 2: estudiante
$ws2 = $wb.Worksheets.Item("estudiante")
$ws2.Range("C1").Value = 10
$ws2.Range("E1").Value = 50
$ws2.Range("F1").Value = 100
$ws2.Range("G1").Value = 30
$ws2.Range("E2").Value = 50
$ws2.Range("G2").Value = ""
$ws2.Range("C3").Value = ""
$ws2.Range("D3").Value = ""
$ws2.Range("E3").Value = ""
$ws2.Range("F3").Value = 80
$ws2.Range("G3").Value = 80
